$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 1.6.3 release renamed the "ProtoBuffMessageType" identifier used in the
# KAFKA protobuf validation step to "ProtobufType".
$ws.Range("K3").Value = "ProtobufType"

# Reflect the updated selection left behind in the saved worksheet view.
$ws.Range("K3").Select()
